$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# New text content for the three "goal" bullets (ilvl=1 sub-bullets under
# "The planned goals I put together for this task were:").
# The three paragraphs get their text content effectively rotated:
#   old para 1 (weapons)  -> becomes the "hunger" paragraph
#   old para 2 (aiming)   -> becomes the "weapons" paragraph (edited)
#   old para 3 (hunger)   -> becomes the "aiming" paragraph (edited)
# ---------------------------------------------------------------------------

$GoalHunger = "The agent will have to contend with hunger, which will be increased with each step it takes, and each shot it fires. When its hunger reaches a threshold, it must consider whether to continue its current actions, or return to the food station and eat if it would not compromise its current attack against the target. When its hunger reaches a high enough level that if it continued its current actions it would not have the stamina to return to a food station to reduce its hunger (based on the maximum possible distance between itself and the food station), it must return to the food station and eat."

$GoalWeapons = "The soldier has two weapons at a time, with a limited number of magazines. When patrolling, they must consider if their ammo remaining, between the two weapons, will theoretically be sufficient to reduce the target’s HP to 0. If it is not, the agent must return to an ammunition station to get new weapons. When attacking, if the soldier runs out of ammo, they must return to the ammo station to exchange their weapons. The choice of weapons should be random, and all weapons, including the ones just handed back, will be available."

$GoalAiming = "When aiming its shot, the agent will consider whether staying with its current weapon or its secondary weapon will be more conducive or sufficient in reducing the target’s HP to 0, in terms of shot impact, speed, rounds per minute. For example, if the target has less than full health, the soldier shouldn’t use two rifle rounds to take out the target if another weapon would also reduce the target’s health to 0 in one shot. "

# New bottom-level ("What I did" style) bullets that replace the "…" placeholder
# and the 4 new ones that follow it.

$Bullet1 = "First, I took the code outlining what each weapon was, as well as the code for pooling projectiles, and created a weapon class that could hold the appropriate variables and be passed from world to soldier and back."
$Bullet2 = "Created points for ammo station and food station"
$Bullet3 = "Added logic for considering when to swap to the next weapon, and to go to the ammo station if attacking and out of ammo, or patrolling and ammo would be insufficient to kill the target when encountered. "
$Bullet4 = "Added logic for handling explosive weapons, for not firing one if too close to the target, and for avoiding existing explosive projectiles so as not to get caught in the blast radius."
$Bullet5 = "Added logic for checking if the soldier will starve if they don’t go and get food now, and for sending them to the food station to satiate its hunger."

# ---------------------------------------------------------------------------
# Locate the three goal paragraphs (ilvl=1, under "Tasks Undertaken") by
# their current (pre-edit) text, so the script does not depend on fixed
# paragraph indices.
# ---------------------------------------------------------------------------

$weaponsOldStart = "The agent has two weapons at a time"
$aimingOldStart = "When aiming its shot"
$hungerOldStart = "The agent will have to contend with hunger"

$pWeapons = $null
$pAiming = $null
$pHunger = $null
$pEllipsis = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith($weaponsOldStart)) {
        $pWeapons = $p
    } elseif ($t.StartsWith($aimingOldStart)) {
        $pAiming = $p
    } elseif ($t.StartsWith($hungerOldStart)) {
        $pHunger = $p
    } elseif ($t.StartsWith([string][char]8230)) {
        $pEllipsis = $p
    }
}

# Replace paragraph text (leave the trailing paragraph mark untouched).
function Set-ParaText($para, [string]$newText) {
    $r = $para.Range
    $target = $d.Range($r.Start, $r.End - 1)
    $target.Text = $newText
}

Set-ParaText $pWeapons $GoalHunger
Set-ParaText $pAiming $GoalWeapons
Set-ParaText $pHunger $GoalAiming

# ---------------------------------------------------------------------------
# Replace the "…" placeholder bullet with real content, then append four
# more bullets at the same (ilvl=0) list level.
# ---------------------------------------------------------------------------

Set-ParaText $pEllipsis $Bullet1

$rEnd = $pEllipsis.Range
$rEnd.Collapse(0)
$rEnd.InsertParagraphAfter()

# Re-find the paragraph following the (now rewritten) ellipsis paragraph.
$pAfterEllipsis = $pEllipsis.Next()
$pAfterEllipsis.Range.Text = $Bullet2

$r2 = $pAfterEllipsis.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$p3 = $pAfterEllipsis.Next()
$p3.Range.Text = $Bullet3

$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$p4 = $p3.Next()
$p4.Range.Text = $Bullet4

$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertParagraphAfter()
$p5 = $p4.Next()
$p5.Range.Text = $Bullet5

# ---------------------------------------------------------------------------
# Drop the stale lastRenderedPageBreak hint on the "T: " bullet (pagination
# shifted because of the extra content added above), and let it reappear at
# the start of the new "Added logic for considering..." bullet instead.
# ---------------------------------------------------------------------------

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("T: ")) {
        $r = $p.Range
        $target = $d.Range($r.Start, $r.End - 1)
        $full = $target.Text
        $target.Delete()
        $ins = $d.Range($r.Start, $r.Start)
        $ins.InsertAfter($full)
        break
    }
}

Write-Host "Done"
